$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 375
$ws.Range("I12").Value = 375
$ws.Range("K12").Value = 375
$ws.Range("M12").Value = -205

$ws.Range("H18").Value = 422.25
$ws.Range("I18").Value = 422.25
$ws.Range("K18").Value = 422.25
$ws.Range("M18").Value = -138.25

$ws.Range("H33").Value = 312.5
$ws.Range("I33").Value = 312.5
$ws.Range("K33").Value = 312.5
$ws.Range("M33").Value = -83.5

$ws.Range("H51").Value = 4887.125
$ws.Range("I51").Value = 3200
$ws.Range("J51").Value = 4999.6
$ws.Range("K51").Value = 3200
$ws.Range("L51").Value = 4999.6
$ws.Range("M51").Value = -2716
$ws.Range("N51").Value = -5967.6

$ws.Range("H64").Value = 9649.916999999999
$ws.Range("I64").Value = 9066.5
$ws.Range("K64").Value = 9066.5
$ws.Range("M64").Value = -8818.5

$ws.Range("H67").Value = 9649.916999999999
$ws.Range("I67").Value = 9066.5
$ws.Range("K67").Value = 9066.5
$ws.Range("M67").Value = -8208.5

$ws.Range("H74").Value = 5030.6
$ws.Range("J74").Value = 6000
$ws.Range("L74").Value = 6000
$ws.Range("N74").Value = -7872

$ws.Range("H77").Value = 5030.6
$ws.Range("J77").Value = 6000
$ws.Range("L77").Value = 30000
$ws.Range("N77").Value = -39360

$ws.Range("H100").Value = 1748.5294
$ws.Range("I100").Value = 942
$ws.Range("K100").Value = 942
$ws.Range("M100").Value = -401

$ws.Range("H111").Value = 4556.75
$ws.Range("I111").Value = 4739.5293
$ws.Range("J111").Value = 4112.857
$ws.Range("K111").Value = 14218.5879
$ws.Range("L111").Value = 12338.571
$ws.Range("M111").Value = -11151.5879
$ws.Range("N111").Value = -18472.571

$ws.Range("H112").Value = 5557603
$ws.Range("J112").Value = 5557603
$ws.Range("L112").Value = 16672809
$ws.Range("N112").Value = -16675025

$ws.Range("H116").Value = 5058.5713
$ws.Range("J116").Value = 4966.6665
$ws.Range("L116").Value = 4966.6665
$ws.Range("N116").Value = -11850.6665

$ws.Range("H138").Value = 9525689
$ws.Range("I138").Value = 1363
$ws.Range("J138").Value = 15875240
$ws.Range("K138").Value = 4089
$ws.Range("L138").Value = 47625720
$ws.Range("M138").Value = 1051
$ws.Range("N138").Value = -47636000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1753
$ws.Range("I2").Value = 1499.5
$ws.Range("K2").Value = 1499.5
$ws.Range("M2").Value = -1386.5

$ws.Range("H45").Value = 2097.6
$ws.Range("I45").Value = 1868.5
$ws.Range("K45").Value = 1868.5
$ws.Range("M45").Value = -1491.5

$ws.Range("H61").Value = 43483028
$ws.Range("I61").Value = 58827012
$ws.Range("K61").Value = 58827012
$ws.Range("M61").Value = -58826800

$ws.Range("H74").Value = 45507460
$ws.Range("I74").Value = 52692076
$ws.Range("K74").Value = 52692076
$ws.Range("M74").Value = -52691202

$ws.Range("H77").Value = 45507460
$ws.Range("I77").Value = 52692076
$ws.Range("K77").Value = 263460380
$ws.Range("M77").Value = -263456012

$ws.Range("H110").Value = 16270.73
$ws.Range("I110").Value = 19768.2
$ws.Range("J110").Value = 4612.5
$ws.Range("K110").Value = 19768.2
$ws.Range("L110").Value = 4612.5
$ws.Range("M110").Value = -17723.2
$ws.Range("N110").Value = -8702.5

$ws.Range("H116").Value = 1753
$ws.Range("I116").Value = 1499.5
$ws.Range("K116").Value = 1499.5
$ws.Range("M116").Value = 794.5

$ws.Range("H122").Value = 2890.1428
$ws.Range("I122").Value = 1763.1428
$ws.Range("J122").Value = 4580.643
$ws.Range("K122").Value = 5289.428400000001
$ws.Range("L122").Value = 13741.929
$ws.Range("M122").Value = -2839.428400000001
$ws.Range("N122").Value = -18641.929

$ws.Range("H132").Value = 33336914
$ws.Range("I132").Value = 3720.2964
$ws.Range("K132").Value = 11160.8892
$ws.Range("M132").Value = -8630.889200000001

$ws.Range("H136").Value = 43483028
$ws.Range("I136").Value = 58827012
$ws.Range("K136").Value = 176481036
$ws.Range("M136").Value = -176478486

$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1753
$ws.Range("I3").Value = 1499.5
$ws.Range("K3").Value = 1499.5
$ws.Range("M3").Value = -1385.5

$ws.Range("H22").Value = 233.33333
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -646

$ws.Range("H74").Value = 47955.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 47955.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 47955.4
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -49827.4

$ws.Range("H77").Value = 47955.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 47955.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 143866.2
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -153226.2

$ws.Range("H86").Value = 17134.5
$ws.Range("I86").Value = 18166.834
$ws.Range("K86").Value = 18166.834
$ws.Range("M86").Value = -17043.834

$ws.Range("H89").Value = 17134.5
$ws.Range("I89").Value = 18166.834
$ws.Range("K89").Value = 90834.17
$ws.Range("M89").Value = -85218.17

$ws.Range("H134").Value = 2956.4878
$ws.Range("I134").Value = 2824.1316
$ws.Range("K134").Value = 8472.3948
$ws.Range("M134").Value = -5937.3948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1905.4615
$ws.Range("J94").Value = 2055.3333
$ws.Range("L94").Value = 2055.3333
$ws.Range("N94").Value = -2957.3333

$ws.Range("H99").Value = 11211
$ws.Range("I99").Value = 11980
$ws.Range("J99").Value = 10885.654
$ws.Range("K99").Value = 11980
$ws.Range("L99").Value = 10885.654
$ws.Range("M99").Value = -10482
$ws.Range("N99").Value = -13881.654

$ws.Range("H107").Value = 1320.8667
$ws.Range("J107").Value = 2359
$ws.Range("L107").Value = 2359
$ws.Range("N107").Value = -6199

$ws.Range("H126").Value = 11211
$ws.Range("I126").Value = 11980
$ws.Range("J126").Value = 10885.654
$ws.Range("K126").Value = 35940
$ws.Range("L126").Value = 32656.962
$ws.Range("M126").Value = -33470
$ws.Range("N126").Value = -37596.962

$ws.Range("H132").Value = 3667.5217
$ws.Range("I132").Value = 2856.389
$ws.Range("K132").Value = 8569.167000000001
$ws.Range("M132").Value = -6039.167000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 1700
$ws.Range("I110").Value = 1700
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5100
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1010
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3232.9167
$ws.Range("I113").Value = 2399.9333
$ws.Range("J113").Value = 4621.222
$ws.Range("K113").Value = 2399.9333
$ws.Range("L113").Value = 4621.222
$ws.Range("M113").Value = -229.9333000000001
$ws.Range("N113").Value = -8961.222

$ws.Range("H132").Value = 3592
$ws.Range("I132").Value = 3374.6667
$ws.Range("K132").Value = 10124.0001
$ws.Range("M132").Value = -7594.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 1041
$ws.Range("I45").Value = 1041
$ws.Range("K45").Value = 1041
$ws.Range("M45").Value = -634

$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H55").Value = 589.1818
$ws.Range("I55").Value = 272.5
$ws.Range("K55").Value = 272.5
$ws.Range("M55").Value = -99.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 461.55554
$ws.Range("I107").Value = 325.83334
$ws.Range("K107").Value = 977.5000200000001
$ws.Range("M107").Value = 942.4999799999999

$ws.Range("H119").Value = 59333.332
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 59333.332
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 59333.332
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -69009.33199999999

$ws.Range("H132").Value = 3927.1538
$ws.Range("I132").Value = 4004.4473
$ws.Range("K132").Value = 12013.3419
$ws.Range("M132").Value = -9483.341899999999
